$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.551.38"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "1.968.67"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'322.57"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.4795"
$ws.Range("E7").Value = "  -4.16%  "
$ws.Range("D8").Value = "'0.4065"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").Value = "'53.84"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'0.08531"
$ws.Range("E10").Value = "  -8.14%  "
$ws.Range("D11").Value = "'1.064"
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D12").Value = "'22.52"
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("D13").Value = "1.964.38"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "'7.633"
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").Value = "'6.214"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").Value = "'1.010"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "'91.29"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'0.00001079"
$ws.Range("E18").Value = "  -3.21%  "
$ws.Range("D19").Value = "'0.06625"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").Value = "'18.68"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "'5.876"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "28.600.73"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'11.59"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'2.294"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "2.201.24"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'155.73"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "'20.41"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "'5.976"
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("D30").Value = "'2.187"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("D31").Value = "'124.92"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "'0.9935"
$ws.Range("E32").Value = "  -5.77%  "
$ws.Range("D33").Value = "'0.09655"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "'1.469"
$ws.Range("E34").Value = "  -4.42%  "
$ws.Range("D35").Value = "'5.683"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").Value = "'3.689"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'9.179"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").Value = "'0.02348"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "'0.06278"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "'1.259"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").Value = "'0.6266"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").Value = "'11.25"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "'1.009"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'0.1923"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("D45").Value = "'1.352"
$ws.Range("E45").Value = "  +5.41%  "
$ws.Range("D46").Value = "'0.5993"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("D47").Value = "'13.12"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").Value = "'2.084"
$ws.Range("E48").Value = "  -4.99%  "
$ws.Range("D49").Value = "'3.412"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "'0.06841"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'0.00000000309"
$ws.Range("E51").Value = "  -7.03%  "
